$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D11").Value = -7.263000000000001
$ws.Range("A12").Value = -21.513
$ws.Range("D23").Value = -8.301
$ws.Range("A27").Value = -21.8
$ws.Range("D28").Value = -8.334999999999999
$ws.Range("A32").Value = -21.604
$ws.Range("D32").Value = -7.306999999999999
$ws.Range("D34").Value = -7.933
$ws.Range("A36").Value = -20.316
$ws.Range("A38").Value = -19.985
$ws.Range("D42").Value = -8.488
$ws.Range("A46").Value = -21.852
$ws.Range("D49").Value = -8.276
$ws.Range("A54").Value = -21.975
$ws.Range("D54").Value = -7.877000000000001
$ws.Range("A55").Value = -22.184
$ws.Range("A56").Value = -22.038
$ws.Range("A67").Value = -21.536
$ws.Range("A69").Value = -21.503
$ws.Range("A72").Value = -21.689
$ws.Range("D78").Value = -8.208000000000002
$ws.Range("D80").Value = -8.242999999999999
$ws.Range("A83").Value = -21.987
$ws.Range("A86").Value = -22.135
$ws.Range("A91").Value = -20.675
$ws.Range("A93").Value = -21.422
$ws.Range("D97").Value = -7.674000000000001
$ws.Range("A99").Value = -22.142
$ws.Range("D99").Value = -8.280000000000001
$ws.Range("D101").Value = -8.019000000000002
$ws.Range("A104").Value = -21.385
